$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying query result reordered the attribute rows (the sensor
# row-count row is no longer returned), so rewrite rows 2-21 with the
# attribute/type pairs in their new order.

$ws.Cells.Item(2, 1).Value = "lifecycle:transition"
$ws.Cells.Item(2, 2).Value = "str"

$ws.Cells.Item(3, 1).Value = "time:timestamp"
$ws.Cells.Item(3, 2).Value = "datetime"

$ws.Cells.Item(4, 1).Value = "case"
$ws.Cells.Item(4, 2).Value = "str"

$ws.Cells.Item(5, 1).Value = "complete_service_time"
$ws.Cells.Item(5, 2).Value = "str"

$ws.Cells.Item(6, 1).Value = "identifier:id"
$ws.Cells.Item(6, 2).Value = "str"

$ws.Cells.Item(7, 1).Value = "unsatisfied_condition_description"
$ws.Cells.Item(7, 2).Value = "str"

$ws.Cells.Item(8, 1).Value = "human_workstation_green_button_pressed"
$ws.Cells.Item(8, 2).Value = "float"

$ws.Cells.Item(9, 1).Value = "response_status_code"
$ws.Cells.Item(9, 2).Value = "float"

$ws.Cells.Item(10, 1).Value = "concept:name"
$ws.Cells.Item(10, 2).Value = "str"

$ws.Cells.Item(11, 1).Value = "SubProcessID"
$ws.Cells.Item(11, 2).Value = "str"

$ws.Cells.Item(12, 1).Value = "lifecycle:state"
$ws.Cells.Item(12, 2).Value = "str"

$ws.Cells.Item(13, 1).Value = "planned_operation_time"
$ws.Cells.Item(13, 2).Value = "str"

$ws.Cells.Item(14, 1).Value = "process_model_id"
$ws.Cells.Item(14, 2).Value = "str"

$ws.Cells.Item(15, 1).Value = "current_task"
$ws.Cells.Item(15, 2).Value = "str"

$ws.Cells.Item(16, 1).Value = "org:resource"
$ws.Cells.Item(16, 2).Value = "str"

$ws.Cells.Item(17, 1).Value = "parameters"
$ws.Cells.Item(17, 2).Value = "dict"

$ws.Cells.Item(18, 1).Value = "case:concept:name"
$ws.Cells.Item(18, 2).Value = "str"

$ws.Cells.Item(19, 1).Value = "operation_end_time"
$ws.Cells.Item(19, 2).Value = "datetime"

$ws.Cells.Item(20, 1).Value = "requested_service_url"
$ws.Cells.Item(20, 2).Value = "str"

$ws.Cells.Item(21, 1).Value = "event_id"
$ws.Cells.Item(21, 2).Value = "str"
